# multi browser implementation for chrome and firefox and accounting cash
#
# Applies to the "Repayment Schedule" sheet: a bunch of placeholder 0-values
# are cleared out (left blank) now that the schedule is driven dynamically,
# an obsolete "O" column's worth of cells is dropped row by row, and the
# long-unused Q/R columns on the opening row are removed outright. The
# active tab moves from "Transactions" back to "Summary", and the
# Repayment Schedule sheet remembers a fresh cell selection (D8).

$wb = $excel.ActiveWorkbook

$wsRepay = $wb.Worksheets.Item("Repayment Schedule")

# --- Row 2 (opening balance row) ------------------------------------------
$wsRepay.Range("A2").Value = $null
$wsRepay.Range("B2").Value = $null
$wsRepay.Range("D2:F2").ClearContents()
$wsRepay.Range("H2").Value = $null
$wsRepay.Range("J2").Value = $null
$wsRepay.Range("M2:O2").ClearContents()
# P2:R2 are dropped completely (format + content), trimming the row back to
# column O and letting the row height revert to the sheet default.
$wsRepay.Range("P2:R2").Clear()
$wsRepay.Rows.Item(2).EntireRow.AutoFit() | Out-Null

# --- Row 3 (first instalment) ----------------------------------------------
# E3's stray 0 is cleared and picks up an italic / wrapped / vertically
# centred look (previously-unused combination among the sheet's styles).
$e3 = $wsRepay.Range("E3")
$e3.Value = $null
$e3.Font.Italic = $true
$e3.VerticalAlignment = -4108   # xlVAlignCenter
$e3.WrapText = $true
$wsRepay.Range("O3").Clear()

# --- Rows 4-14 (remaining instalments) -------------------------------------
for ($r = 4; $r -le 14; $r++) {
    $wsRepay.Range("D$r`:E$r").ClearContents()
    $wsRepay.Range("O$r").Clear()
}

# Selection on the Repayment Schedule sheet moves to D8.
$wsRepay.Activate() | Out-Null
$wsRepay.Range("D8").Select() | Out-Null

# --- Active tab moves to Summary -------------------------------------------
# Doing this last also flips the Transactions sheet's tabSelected flag off
# and Summary's on, and updates the workbook-level activeTab.
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Activate() | Out-Null
